$d = $word.ActiveDocument

# Update the date line (unique text, safe to use Find/Replace).
$d.Content.Find.Execute("2023-08-20 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-08-21 Monday", 2)

# Update the division problems/answers table cell-by-cell (positional,
# since several values repeat so text-based Find/Replace would be
# ambiguous).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "34÷4=8, 2"
$t.Cell(1, 2).Range.Text = "26÷7=3, 5"
$t.Cell(1, 3).Range.Text = "15÷6=2, 3"
$t.Cell(1, 4).Range.Text = "15÷4=3, 3"
$t.Cell(1, 5).Range.Text = "53÷7=7, 4"

$t.Cell(5, 1).Range.Text = "30÷5=6, 0"
$t.Cell(5, 2).Range.Text = "11÷4=2, 3"
$t.Cell(5, 3).Range.Text = "14÷5=2, 4"
$t.Cell(5, 4).Range.Text = "96÷5=19, 1"
$t.Cell(5, 5).Range.Text = "25÷9=2, 7"

$t.Cell(9, 1).Range.Text = "49÷4=12, 1"
$t.Cell(9, 2).Range.Text = "71÷8=8, 7"
$t.Cell(9, 3).Range.Text = "18÷7=2, 4"
$t.Cell(9, 4).Range.Text = "18÷3=6, 0"
$t.Cell(9, 5).Range.Text = "36÷2=18, 0"

$t.Cell(13, 1).Range.Text = "72÷6=12, 0"
$t.Cell(13, 2).Range.Text = "11÷5=2, 1"
$t.Cell(13, 3).Range.Text = "88÷2=44, 0"
$t.Cell(13, 4).Range.Text = "61÷2=30, 1"
$t.Cell(13, 5).Range.Text = "44÷9=4, 8"

$t.Cell(17, 1).Range.Text = "75÷6=12, 3"
$t.Cell(17, 2).Range.Text = "46÷9=5, 1"
$t.Cell(17, 3).Range.Text = "44÷9=4, 8"
$t.Cell(17, 4).Range.Text = "30÷4=7, 2"
$t.Cell(17, 5).Range.Text = "45÷4=11, 1"

Write-Output "done"
